$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.519.19'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '2.250.77'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.06'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.11'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.15'
$ws.Range("E10").Value = '  +2.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").Value = '2.312.29'
$ws.Range("E14").Value = '  +2.36%  '
$ws.Range("E15").Value = '  +3.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.69'
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("D17").Value = '44.217.18'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").Value = '0.0₃0965'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.41'
$ws.Range("E19").Value = '  +5.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.17'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.99'
$ws.Range("E21").Value = '  +2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.51'
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("E23").Value = '  +4.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.00'
$ws.Range("E24").Value = '  +4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  +5.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.87'
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.84'
$ws.Range("E28").Value = '  +5.15%  '
$ws.Range("E29").Value = '  +2.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.10'
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.52'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("E37").Value = '  +2.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.44'
$ws.Range("E38").Value = '  +4.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.53'
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0304'
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").Value = '1.755.89'
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("E44").Value = '  +5.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '80.85'
$ws.Range("E45").Value = '  -3.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.11'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.06'
$ws.Range("E47").Value = '  +4.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.67'
$ws.Range("E48").Value = '  +4.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.18'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.88'
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("E51").Value = '  +5.89%  '
